# Assignment 3 / Part 2 - User Evaluation : add a 10th task (T10) column
# to both data tables (CLICKS @ rows 2-7, TIME @ rows 10-15), update the
# values that shifted when the new task was inserted, refresh the
# AVG TIME / AVG TASK TIME formulas to include the new column, and point
# both line charts at the widened ranges.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. New header label "T10" in column K, mirroring the existing T1..T9
#    headers in row 2 (CLICKS table) and row 10 (TIME table).
# ---------------------------------------------------------------------
$ws.Range("K2").Value = "T10"
$ws.Range("K2").Font.Bold = $true

$ws.Range("K10").Value = "T10"
$ws.Range("K10").Font.Bold = $true

# ---------------------------------------------------------------------
# 2. CLICKS table (rows 3-7, columns B..K) - updated values
# ---------------------------------------------------------------------
$clicks = @{
    3 = @(5,1,2,1,3,1,1,2,2,2)
    4 = @(2,1,1,1,3,1,1,2,2,2)
    5 = @(2,1,1,1,3,1,1,2,2,2)
    6 = @(2,1,2,1,3,1,1,2,2,2)
    7 = @(2,1,1,1,3,1,1,2,2,2)
}
$cols = @("B","C","D","E","F","G","H","I","J","K")

foreach ($row in $clicks.Keys) {
    $vals = $clicks[$row]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range("$($cols[$i])$row").Value = $vals[$i]
    }
}

# ---------------------------------------------------------------------
# 3. TIME (seconds) table (rows 11-15, columns B..K) - updated values
# ---------------------------------------------------------------------
$times = @{
    11 = @(35,11,12,4,11,5,4,14,8,7)
    12 = @(21,6,4,2,9,4,3,9,6,6)
    13 = @(7,3,3,2,6,2,1,5,3,2)
    14 = @(20,15,9,4,8,4,3,10,6,8)
    15 = @(9,6,5,3,9,4,2,6,4,4)
}

foreach ($row in $times.Keys) {
    $vals = $times[$row]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range("$($cols[$i])$row").Value = $vals[$i]
    }
}

# ---------------------------------------------------------------------
# 4. AVG TIME row (16) - extend the shared "(SUM(col11:col15)/5)"
#    formula from column J out to the new column K.
# ---------------------------------------------------------------------
$ws.Range("K16").Formula = "=(SUM(K11:K15)/5)"

# AVG TASK TIME (B17) - widen the average to cover the 10 tasks.
$ws.Range("B17").Formula = "=SUM(B16:K16)/10"

# ---------------------------------------------------------------------
# 5. Point both line charts at the widened source ranges (B..K instead
#    of B..J) so the cached chart data/categories include T10.
# ---------------------------------------------------------------------
$clicksChart = $ws.ChartObjects().Item(1).Chart
$clicksSeries = $clicksChart.SeriesCollection()
for ($i = 1; $i -le $clicksSeries.Count; $i++) {
    $s = $clicksSeries.Item($i)
    $dataRow = 2 + $i
    $s.Values = $ws.Range("B" + $dataRow + ":K" + $dataRow)
    $s.XValues = $ws.Range("B10:K10")
}

$timeChart = $ws.ChartObjects().Item(2).Chart
$timeSeries = $timeChart.SeriesCollection()
for ($i = 1; $i -le $timeSeries.Count; $i++) {
    $s = $timeSeries.Item($i)
    $dataRow = 10 + $i
    $s.Values = $ws.Range("B" + $dataRow + ":K" + $dataRow)
    $s.XValues = $ws.Range("B10:K10")
}

# ---------------------------------------------------------------------
# 6. Selection / scroll position left by the author after editing.
# ---------------------------------------------------------------------
$ws.Range("M16").Select()
